# Apply "provider, banners validation - OK" changes:
#  - Replace the "fixture.channel.count()" note (used in sheet "Каналы", cells G5 & G8)
#    with "добавление каналов через апи"
#  - Switch the active/selected sheet from "разное" to "Каналы"
#  - Update the selection on "Каналы" to the single cell G5

$wb = $excel.ActiveWorkbook

$channels = $wb.Worksheets.Item("Каналы")

# Update the two cells that referenced the old shared string so it becomes
# the new note text (the now-unused old string is dropped on save).
$channels.Range("G5").Value = "добавление каналов через апи"
$channels.Range("G8").Value = "добавление каналов через апи"

# Make "Каналы" the active sheet (this flips tabSelected/activeTab)
$channels.Activate()

# Update the selection shown on the "Каналы" sheet
$channels.Range("G5").Select()
